$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new "2022" column (N) is being added, mirroring the existing "2021"
# column (M): same formatting/styles, new data values.
$ws.Range("M4:M13").Copy()
$ws.Range("N4:N13").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("N4").Value = 2022
$ws.Range("N5").Value = 4.3
$ws.Range("N6").Value = 5.1
$ws.Range("N7").Value = 3.1
$ws.Range("N8").Value = 2.9
$ws.Range("N9").Value = 3.4
$ws.Range("N10").Value = 2.3
$ws.Range("N11").Value = 92.8
$ws.Range("N12").Value = 91.6
$ws.Range("N13").Value = 94.6

$ws.Range("N15").Select()
